$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.Value = "'29.981.17"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.47%  "

# Row 3
$cell = $ws.Range("D3")
$cell.Value = "'1.880.69"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -1.37%  "

# Row 4
$cell = $ws.Range("D4")
$cell.Value = "'1.000"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$cell = $ws.Range("D5")
$cell.Value = "'242.63"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -4.11%  "

# Row 6
$cell = $ws.Range("D6")
$cell.Value = "'0.9998"
$cell.Style = "Normal"

# Row 7
$cell = $ws.Range("D7")
$cell.Value = "'0.4909"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -3.51%  "

# Row 8
$cell = $ws.Range("D8")
$cell.Value = "'0.2939"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -3.30%  "

# Row 9
$cell = $ws.Range("D9")
$cell.Value = "'0.06616"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -3.02%  "

# Row 10
$cell = $ws.Range("D10")
$cell.Value = "'1.881.85"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -1.35%  "

# Row 11
$cell = $ws.Range("D11")
$cell.Value = "'16.68"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -3.79%  "

# Row 12
$cell = $ws.Range("D12")
$cell.Value = "'0.07170"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -2.13%  "

# Row 13
$cell = $ws.Range("D13")
$cell.Value = "'0.6663"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -3.73%  "

# Row 14
$cell = $ws.Range("D14")
$cell.Value = "'86.64"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -0.36%  "

# Row 15
$cell = $ws.Range("D15")
$cell.Value = "'4.876"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.85%  "

# Row 16
$cell = $ws.Range("D16")
$cell.Value = "'29.961.84"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.57%  "

# Row 17
$cell = $ws.Range("D17")
$cell.Value = "'0.000007807"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -5.56%  "

# Row 18
$cell = $ws.Range("D18")
$cell.Value = "'0.9997"
$cell.Style = "Normal"

# Row 19
$ws.Range("E19").Value = "  -2.21%  "

# Row 20
$cell = $ws.Range("D20")
$cell.Value = "'2.121.03"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.49%  "

# Row 21
$cell = $ws.Range("D21")
$cell.Value = "'1.000"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.26%  "

# Row 22
$cell = $ws.Range("D22")
$cell.Value = "'4.779"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.95%  "

# Row 23
$cell = $ws.Range("D23")
$cell.Value = "'5.841"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +1.66%  "

# Row 24
$cell = $ws.Range("D24")
$cell.Value = "'9.089"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -2.28%  "

# Row 25
$cell = $ws.Range("D25")
$cell.Value = "'150.51"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.96%  "

# Row 26
$cell = $ws.Range("D26")
$cell.Value = "'141.14"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +4.45%  "

# Row 27
$cell = $ws.Range("D27")
$cell.Value = "'16.98"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.84%  "

# Row 28
$cell = $ws.Range("D28")
$cell.Value = "'1.902"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -5.12%  "

# Row 29
$cell = $ws.Range("D29")
$cell.Value = "'1.392"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -0.71%  "

# Row 30
$cell = $ws.Range("D30")
$cell.Value = "'4.195"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -2.15%  "

# Row 31
$cell = $ws.Range("D31")
$cell.Value = "'0.08740"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -1.29%  "

# Row 32
$cell = $ws.Range("D32")
$cell.Value = "'3.981"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.66%  "

# Row 33
$cell = $ws.Range("D33")
$cell.Value = "'0.05021"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.72%  "

# Row 34
$cell = $ws.Range("D34")
$cell.Value = "'0.7168"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -1.19%  "

# Row 35
$cell = $ws.Range("D35")
$cell.Value = "'1.111"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -2.92%  "

# Row 36
$cell = $ws.Range("D36")
$cell.Value = "'2.670"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.59%  "

# Row 37
$cell = $ws.Range("D37")
$cell.Value = "'0.01798"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +5.99%  "

# Row 38
$ws.Range("E38").Value = "  -4.46%  "

# Row 39
$cell = $ws.Range("D39")
$cell.Value = "'2.156"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -5.24%  "

# Row 40
$cell = $ws.Range("D40")
$cell.Value = "'0.9388"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -2.54%  "

# Row 41
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$cell = $ws.Range("D41")
$cell.Value = "'0.9991"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.01%  "

# Row 42
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$cell = $ws.Range("D42")
$cell.Value = "'103.76"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -1.24%  "

# Row 43
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$cell = $ws.Range("D43")
$cell.Value = "'0.4228"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -2.09%  "

# Row 44
$cell = $ws.Range("D44")
$cell.Value = "'5.734"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -6.84%  "

# Row 45
$cell = $ws.Range("D45")
$cell.Value = "'7.331"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -4.21%  "

# Row 46
$ws.Range("E46").Value = "  -1.00%  "

# Row 47
$cell = $ws.Range("D47")
$cell.Value = "'0.05700"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -0.79%  "

# Row 48
$cell = $ws.Range("D48")
$cell.Value = "'32.61"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -1.83%  "

# Row 49
$cell = $ws.Range("D49")
$cell.Value = "'8.274"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -2.19%  "

# Row 50
$cell = $ws.Range("D50")
$cell.Value = "'0.3749"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -2.17%  "

# Row 51
$cell = $ws.Range("D51")
$cell.Value = "'55.98"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -1.61%  "
